# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.607.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "'1.923.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'246.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4748"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "'0.2898"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "'0.06848"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.10%  "
$ws.Range("D10").Value = "'105.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").Value = "'18.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.71%  "
$ws.Range("D12").Value = "'1.922.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'0.07688"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "'5.346"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.32%  "
$ws.Range("D15").Value = "'0.6693"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "'290.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("D17").Value = "'30.623.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "'0.000007620"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'5.573"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.92%  "
$ws.Range("D21").Value = "'12.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "'2.177.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("D25").Value = "'9.537"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.47%  "
$ws.Range("D26").Value = "'167.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").Value = "'21.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.80%  "
$ws.Range("D28").Value = "'2.117"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.00%  "
$ws.Range("D29").Value = "'0.1072"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("D30").Value = "'1.403"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("D31").Value = "'4.176"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.10%  "
$ws.Range("D32").Value = "'4.047"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("D33").Value = "'0.05032"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'0.7307"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "'1.144"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "'0.02067"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.37%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'2.735"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.16%  "
$ws.Range("D39").Value = "'2.687"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").Value = "'111.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.39%  "
$ws.Range("D41").Value = "'2.050"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'0.8727"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "'0.4401"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.38%  "
$ws.Range("D44").Value = "'5.935"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("E45").Value = "  +0.12%  "
$ws.Range("D46").Value = "'67.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("D47").Value = "'7.298"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.371"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("B49").Value = "BitcoinSV"
$ws.Range("C49").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D49").Value = "'48.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.62%  "
$ws.Range("D50").Value = "'0.1246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("D51").Value = "'34.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "
